$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptos list: price and volume(1h) changes, plus two coin-row reorderings
# (rows 39/40 VeChain<->Maker, rows 49/50 Frax<->EnergySwap).
# NumberFormat is forced to text ("@") before each write so numeric-looking
# values (e.g. "224.02", "1.000", "0.06664") are stored verbatim as strings
# and not coerced into floating point numbers by Excel.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '27.334.30'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -0.67%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.711.05'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -0.69%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '224.02'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -0.93%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5293'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -0.90%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.004'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -0.09%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2664'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -0.21%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06664'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +1.12%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '20.86'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -3.86%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07683'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '4.496'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -2.66%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.946.43'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -0.70%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '1.712.68'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -0.65%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.5831'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +0.12%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0₅8212'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -1.05%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '68.02'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +0.19%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '27.345.77'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -0.64%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '221.56'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +0.71%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '1.007'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +0.20%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '4.632'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -2.11%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '10.45'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -1.87%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '6.004'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -1.34%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.004'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -0.07%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '144.57'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -1.87%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.690'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -2.66%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.1206'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -2.38%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.229'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -2.41%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '16.23'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -1.98%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.05330'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -4.02%  '
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -1.05%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.466'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -2.79%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.436'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -0.26%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.635'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -1.74%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.877'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +0.50%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.9515'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -1.41%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.392'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -1.18%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.5842'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -1.99%  '
$ws.Range('B39').NumberFormat = '@'
$ws.Range('B39').Value = 'Maker'
$ws.Range('C39').NumberFormat = '@'
$ws.Range('C39').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.148.29'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +8.71%  '
$ws.Range('B40').NumberFormat = '@'
$ws.Range('B40').Value = 'VeChain'
$ws.Range('C40').NumberFormat = '@'
$ws.Range('C40').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.01636'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -1.12%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '5.807'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -1.83%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.004'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -0.11%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.8401'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -1.62%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '101.27'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -0.32%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.853.60'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -0.72%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0₈115'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -0.48%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '57.72'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -2.10%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.4544'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +2.33%  '
$ws.Range('B49').NumberFormat = '@'
$ws.Range('B49').Value = 'EnergySwap'
$ws.Range('C49').NumberFormat = '@'
$ws.Range('C49').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '8.107'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -1.41%  '
$ws.Range('B50').NumberFormat = '@'
$ws.Range('B50').Value = 'Frax'
$ws.Range('C50').NumberFormat = '@'
$ws.Range('C50').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.000'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -0.51%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.05225'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -0.38%  '
